$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_21_9_21"
$ws.Cells.Item(2, 2).Value = 0.957864779800686
$ws.Cells.Item(2, 3).Value = 0.8153829851795996
$ws.Cells.Item(2, 4).Value = 0.8032460883659864
$ws.Cells.Item(2, 5).Value = 0.725615641684737
$ws.Cells.Item(2, 6).Value = 0.779339507905594
$ws.Cells.Item(2, 7).Value = 0.2817584114876626
$ws.Cells.Item(2, 8).Value = 1.234534828187205
$ws.Cells.Item(2, 9).Value = 0.4825631835328374
$ws.Cells.Item(2, 10).Value = 0.3446703797135212
$ws.Cells.Item(2, 11).Value = 0.4136166161637796
$ws.Cells.Item(2, 12).Value = 1.299265393789203
$ws.Cells.Item(2, 13).Value = 0.5308092044112108
$ws.Cells.Item(2, 14).Value = 1.017741145347079
$ws.Cells.Item(2, 15).Value = 0.5534068511695283
$ws.Cells.Item(2, 16).Value = 164.5334105441415
$ws.Cells.Item(2, 17).Value = 263.2623523584658

$ws.Cells.Item(3, 1).Value = "model_21_9_22"
$ws.Cells.Item(3, 2).Value = 0.9576656128919266
$ws.Cells.Item(3, 3).Value = 0.8153815619055185
$ws.Cells.Item(3, 4).Value = 0.7932629099548328
$ws.Cells.Item(3, 5).Value = 0.7236800184718275
$ws.Cells.Item(3, 6).Value = 0.7721595927321212
$ws.Cells.Item(3, 7).Value = 0.2830902415236153
$ws.Cells.Item(3, 8).Value = 1.234544345627539
$ws.Cells.Item(3, 9).Value = 0.5070481572538363
$ws.Cells.Item(3, 10).Value = 0.3471018302228438
$ws.Cells.Item(3, 11).Value = 0.42707499373834
$ws.Cells.Item(3, 12).Value = 1.303119301048126
$ws.Cells.Item(3, 13).Value = 0.5320622534286897
$ws.Cells.Item(3, 14).Value = 1.017825005098136
$ws.Cells.Item(3, 15).Value = 0.5547132450778502
$ws.Cells.Item(3, 16).Value = 164.5239791150537
$ws.Cells.Item(3, 17).Value = 263.252920929378

$ws.Cells.Item(4, 1).Value = "model_21_9_20"
$ws.Cells.Item(4, 2).Value = 0.9580392723886167
$ws.Cells.Item(4, 3).Value = 0.8153617007197469
$ws.Cells.Item(4, 4).Value = 0.8130290824384099
$ws.Cells.Item(4, 5).Value = 0.7274210240970035
$ws.Cells.Item(4, 6).Value = 0.7863446404971989
$ws.Cells.Item(4, 7).Value = 0.280591578748705
$ws.Cells.Item(4, 8).Value = 1.234677157468247
$ws.Cells.Item(4, 9).Value = 0.4585691865400204
$ws.Cells.Item(4, 10).Value = 0.3424025323574083
$ws.Cells.Item(4, 11).Value = 0.4004858594487144
$ws.Cells.Item(4, 12).Value = 1.295406607520545
$ws.Cells.Item(4, 13).Value = 0.5297089566438394
$ws.Cells.Item(4, 14).Value = 1.01766767478374
$ws.Cells.Item(4, 15).Value = 0.552259763576873
$ws.Cells.Item(4, 16).Value = 164.5417102467021
$ws.Cells.Item(4, 17).Value = 263.2706520610264

$ws.Cells.Item(5, 1).Value = "model_21_9_23"
$ws.Cells.Item(5, 2).Value = 0.9574426268009075
$ws.Cells.Item(5, 3).Value = 0.8153579815734057
$ws.Cells.Item(5, 4).Value = 0.7830886990427304
$ws.Cells.Item(5, 5).Value = 0.7216191161329948
$ws.Cells.Item(5, 6).Value = 0.7648127165698435
$ws.Cells.Item(5, 7).Value = 0.2845813505410162
$ws.Cells.Item(5, 8).Value = 1.234702027416955
$ws.Cells.Item(5, 9).Value = 0.5320016616945071
$ws.Cells.Item(5, 10).Value = 0.3496906512330481
$ws.Cells.Item(5, 11).Value = 0.4408463310029897
$ws.Cells.Item(5, 12).Value = 1.306955668694294
$ws.Cells.Item(5, 13).Value = 0.5334616673585987
$ws.Cells.Item(5, 14).Value = 1.017918893978565
$ws.Cells.Item(5, 15).Value = 0.5561722349559416
$ws.Cells.Item(5, 16).Value = 164.5134722482834
$ws.Cells.Item(5, 17).Value = 263.2424140626077

$ws.Cells.Item(6, 1).Value = "model_21_9_19"
$ws.Cells.Item(6, 2).Value = 0.9581884007747294
$ws.Cells.Item(6, 3).Value = 0.8153172169107126
$ws.Cells.Item(6, 4).Value = 0.8226049270028466
$ws.Cells.Item(6, 5).Value = 0.7290949928501114
$ws.Cells.Item(6, 6).Value = 0.7931703388760153
$ws.Cells.Item(6, 7).Value = 0.2795943565440963
$ws.Cells.Item(6, 8).Value = 1.234974620904088
$ws.Cells.Item(6, 9).Value = 0.4350832492102163
$ws.Cells.Item(6, 10).Value = 0.340299761451279
$ws.Cells.Item(6, 11).Value = 0.3876914428333795
$ws.Cells.Item(6, 12).Value = 1.291544330439193
$ws.Cells.Item(6, 13).Value = 0.5287668262515116
$ws.Cells.Item(6, 14).Value = 1.017604883884324
$ws.Cells.Item(6, 15).Value = 0.551277524743265
$ws.Cells.Item(6, 16).Value = 164.5488309057203
$ws.Cells.Item(6, 17).Value = 263.2777727200446

$ws.Cells.Item(7, 1).Value = "model_21_9_24"
$ws.Cells.Item(7, 2).Value = 0.9571964060390739
$ws.Cells.Item(7, 3).Value = 0.8153126682410371
$ws.Cells.Item(7, 4).Value = 0.7727302594457981
$ws.Cells.Item(7, 5).Value = 0.7194324869143602
$ws.Cells.Item(7, 6).Value = 0.7573033716913908
$ws.Cells.Item(7, 7).Value = 0.2862278299091404
$ws.Cells.Item(7, 8).Value = 1.235005037879152
$ws.Cells.Item(7, 9).Value = 0.5574070096584464
$ws.Cells.Item(7, 10).Value = 0.3524374051941959
$ws.Cells.Item(7, 11).Value = 0.4549222074263212
$ws.Cells.Item(7, 12).Value = 1.310780468161841
$ws.Cells.Item(7, 13).Value = 0.5350026447683603
$ws.Cells.Item(7, 14).Value = 1.018022565878285
$ws.Cells.Item(7, 15).Value = 0.5577788149642999
$ws.Cells.Item(7, 16).Value = 164.5019343544243
$ws.Cells.Item(7, 17).Value = 263.2308761687486

$ws.Cells.Item(8, 1).Value = "model_21_9_18"
$ws.Cells.Item(8, 2).Value = 0.9583115053809672
$ws.Cells.Item(8, 3).Value = 0.8152490509128963
$ws.Cells.Item(8, 4).Value = 0.8319650635301257
$ws.Cells.Item(8, 5).Value = 0.7306354332423843
$ws.Cells.Item(8, 6).Value = 0.7998101069560549
$ws.Cells.Item(8, 7).Value = 0.2787711554753398
$ws.Cells.Item(8, 8).Value = 1.235430447245361
$ws.Cells.Item(8, 9).Value = 0.4121263623895485
$ws.Cells.Item(8, 10).Value = 0.3383647233966655
$ws.Cells.Item(8, 11).Value = 0.375245542893107
$ws.Cells.Item(8, 12).Value = 1.287682410110113
$ws.Cells.Item(8, 13).Value = 0.5279878364842696
$ws.Cells.Item(8, 14).Value = 1.017553050365908
$ws.Cells.Item(8, 15).Value = 0.5504653717689003
$ws.Cells.Item(8, 16).Value = 164.5547281300114
$ws.Cells.Item(8, 17).Value = 263.2836699443356

$ws.Cells.Item(9, 1).Value = "model_21_9_17"
$ws.Cells.Item(9, 2).Value = 0.9584078482253184
$ws.Cells.Item(9, 3).Value = 0.8151567488808869
$ws.Cells.Item(9, 4).Value = 0.8411031145833545
$ws.Cells.Item(9, 5).Value = 0.7320387192019457
$ws.Cells.Item(9, 6).Value = 0.8062587358486955
$ws.Cells.Item(9, 7).Value = 0.2781269104315451
$ws.Cells.Item(9, 8).Value = 1.23604767135842
$ws.Cells.Item(9, 9).Value = 0.3897141675268904
$ws.Cells.Item(9, 10).Value = 0.3366019731163709
$ws.Cells.Item(9, 11).Value = 0.3631579234187117
$ws.Cells.Item(9, 12).Value = 1.283822610884851
$ws.Cells.Item(9, 13).Value = 0.5273773890029275
$ws.Cells.Item(9, 14).Value = 1.017512484957761
$ws.Cells.Item(9, 15).Value = 0.549828936274477
$ws.Cells.Item(9, 16).Value = 164.5593555142357
$ws.Cells.Item(9, 17).Value = 263.28829732856

$ws.Cells.Item(10, 1).Value = "model_21_9_16"
$ws.Cells.Item(10, 2).Value = 0.9584768008776193
$ws.Cells.Item(10, 3).Value = 0.815039847054723
$ws.Cells.Item(10, 4).Value = 0.8500117855192612
$ws.Cells.Item(10, 5).Value = 0.7333033879722812
$ws.Cells.Item(10, 6).Value = 0.8125106946249641
$ws.Cells.Item(10, 7).Value = 0.2776658237281108
$ws.Cells.Item(10, 8).Value = 1.236829394408261
$ws.Cells.Item(10, 9).Value = 0.3678645556326463
$ws.Cells.Item(10, 10).Value = 0.3350133480651478
$ws.Cells.Item(10, 11).Value = 0.3514389518488971
$ws.Cells.Item(10, 12).Value = 1.279974771372608
$ws.Cells.Item(10, 13).Value = 0.526940057054036
$ws.Cells.Item(10, 14).Value = 1.017483452262055
$ws.Cells.Item(10, 15).Value = 0.5493729861991193
$ws.Cells.Item(10, 16).Value = 164.5626739225623
$ws.Cells.Item(10, 17).Value = 263.2916157368865

$ws.Cells.Item(11, 1).Value = "model_21_9_15"
$ws.Cells.Item(11, 2).Value = 0.9585176597689962
$ws.Cells.Item(11, 3).Value = 0.8148978867320622
$ws.Cells.Item(11, 4).Value = 0.8586838463438589
$ws.Cells.Item(11, 5).Value = 0.7344270262559089
$ws.Cells.Item(11, 6).Value = 0.818560606390646
$ws.Cells.Item(11, 7).Value = 0.2773926001333366
$ws.Cells.Item(11, 8).Value = 1.237778683739561
$ws.Cells.Item(11, 9).Value = 0.3465952591568917
$ws.Cells.Item(11, 10).Value = 0.3336018797283354
$ws.Cells.Item(11, 11).Value = 0.3400987068922232
$ws.Cells.Item(11, 12).Value = 1.276135431859574
$ws.Cells.Item(11, 13).Value = 0.5266807383352239
$ws.Cells.Item(11, 14).Value = 1.017466248518317
$ws.Cells.Item(11, 15).Value = 0.5491026277455836
$ws.Cells.Item(11, 16).Value = 164.5646428945126
$ws.Cells.Item(11, 17).Value = 263.2935847088368

$ws.Cells.Item(12, 1).Value = "model_21_9_14"
$ws.Cells.Item(12, 2).Value = 0.9585297319528123
$ws.Cells.Item(12, 3).Value = 0.8147303579826943
$ws.Cells.Item(12, 4).Value = 0.8671127913885486
$ws.Cells.Item(12, 5).Value = 0.7354075701466692
$ws.Cells.Item(12, 6).Value = 0.8244037074943962
$ws.Cells.Item(12, 7).Value = 0.27731187338457
$ws.Cells.Item(12, 8).Value = 1.238898949257983
$ws.Cells.Item(12, 9).Value = 0.3259222340525427
$ws.Cells.Item(12, 10).Value = 0.3323701607002201
$ws.Cells.Item(12, 11).Value = 0.3291461177653848
$ws.Cells.Item(12, 12).Value = 1.272321637884501
$ws.Cells.Item(12, 13).Value = 0.5266040954878437
$ws.Cells.Item(12, 14).Value = 1.017461165493553
$ws.Cells.Item(12, 15).Value = 0.5490227220535179
$ws.Cells.Item(12, 16).Value = 164.5652250188919
$ws.Cells.Item(12, 17).Value = 263.2941668332161

$ws.Cells.Item(13, 1).Value = "model_21_9_13"
$ws.Cells.Item(13, 2).Value = 0.9585124478244249
$ws.Cells.Item(13, 3).Value = 0.8145368368618477
$ws.Cells.Item(13, 4).Value = 0.8752927843393362
$ws.Cells.Item(13, 5).Value = 0.7362423720122417
$ws.Cells.Item(13, 6).Value = 0.830034945118199
$ws.Cells.Item(13, 7).Value = 0.277427452430683
$ws.Cells.Item(13, 8).Value = 1.240193025884167
$ws.Cells.Item(13, 9).Value = 0.3058597945979676
$ws.Cells.Item(13, 10).Value = 0.331321516827956
$ws.Cells.Item(13, 11).Value = 0.3185906557129617
$ws.Cells.Item(13, 12).Value = 1.268522973737352
$ws.Cells.Item(13, 13).Value = 0.5267138240360537
$ws.Cells.Item(13, 14).Value = 1.017468443021295
$ws.Cells.Item(13, 15).Value = 0.5491371219731946
$ws.Cells.Item(13, 16).Value = 164.5643916253147
$ws.Cells.Item(13, 17).Value = 263.293333439639

$ws.Cells.Item(14, 1).Value = "model_21_9_12"
$ws.Cells.Item(14, 2).Value = 0.9584651343691772
$ws.Cells.Item(14, 3).Value = 0.8143169050440902
$ws.Cells.Item(14, 4).Value = 0.8832171175957355
$ws.Cells.Item(14, 5).Value = 0.7369285827581398
$ws.Cells.Item(14, 6).Value = 0.8354492427298483
$ws.Cells.Item(14, 7).Value = 0.2777438377237824
$ws.Cells.Item(14, 8).Value = 1.241663710962205
$ws.Cells.Item(14, 9).Value = 0.2864243920088883
$ws.Cells.Item(14, 10).Value = 0.3304595270272093
$ws.Cells.Item(14, 11).Value = 0.3084418364305507
$ws.Cells.Item(14, 12).Value = 1.264753551583945
$ws.Cells.Item(14, 13).Value = 0.5270140773487767
$ws.Cells.Item(14, 14).Value = 1.017488364476136
$ws.Cells.Item(14, 15).Value = 0.5494501576910502
$ws.Cells.Item(14, 16).Value = 164.5621120742431
$ws.Cells.Item(14, 17).Value = 263.2910538885673

$ws.Cells.Item(15, 1).Value = "model_21_9_11"
$ws.Cells.Item(15, 2).Value = 0.9583871595093046
$ws.Cells.Item(15, 3).Value = 0.8140700902870995
$ws.Cells.Item(15, 4).Value = 0.890880985326262
$ws.Cells.Item(15, 5).Value = 0.7374656922504645
$ws.Cells.Item(15, 6).Value = 0.840643045829695
$ws.Cells.Item(15, 7).Value = 0.2782652559707917
$ws.Cells.Item(15, 8).Value = 1.243314162378678
$ws.Cells.Item(15, 9).Value = 0.2676278131784923
$ws.Cells.Item(15, 10).Value = 0.3297848321072656
$ws.Cells.Item(15, 11).Value = 0.2987063226428789
$ws.Cells.Item(15, 12).Value = 1.261013849064705
$ws.Cells.Item(15, 13).Value = 0.5275085363961343
$ws.Cells.Item(15, 14).Value = 1.017521195996082
$ws.Cells.Item(15, 15).Value = 0.5499656668837251
$ws.Cells.Item(15, 16).Value = 164.5583609240855
$ws.Cells.Item(15, 17).Value = 263.2873027384097

$ws.Cells.Item(16, 1).Value = "model_21_9_10"
$ws.Cells.Item(16, 2).Value = 0.9582779155519461
$ws.Cells.Item(16, 3).Value = 0.8137959869896298
$ws.Cells.Item(16, 4).Value = 0.8982782200985022
$ws.Cells.Item(16, 5).Value = 0.7378500254629965
$ws.Cells.Item(16, 6).Value = 0.8456112851175008
$ws.Cells.Item(16, 7).Value = 0.2789957708166708
$ws.Cells.Item(16, 8).Value = 1.245147092390988
$ws.Cells.Item(16, 9).Value = 0.2494851844938238
$ws.Cells.Item(16, 10).Value = 0.3293020484853665
$ws.Cells.Item(16, 11).Value = 0.2893936164895951
$ws.Cells.Item(16, 12).Value = 1.257315181689437
$ws.Cells.Item(16, 13).Value = 0.5282005024767307
$ws.Cells.Item(16, 14).Value = 1.017567193451812
$ws.Cells.Item(16, 15).Value = 0.5506870913929398
$ws.Cells.Item(16, 16).Value = 164.5531173112777
$ws.Cells.Item(16, 17).Value = 263.2820591256019

$ws.Cells.Item(17, 1).Value = "model_21_9_9"
$ws.Cells.Item(17, 2).Value = 0.9581368958937938
$ws.Cells.Item(17, 3).Value = 0.8134941082706717
$ws.Cells.Item(17, 4).Value = 0.9054040065113339
$ws.Cells.Item(17, 5).Value = 0.7380810594777809
$ws.Cells.Item(17, 6).Value = 0.8503505693530051
$ws.Cells.Item(17, 7).Value = 0.279938769920071
$ws.Cells.Item(17, 8).Value = 1.247165756774685
$ws.Cells.Item(17, 9).Value = 0.2320083163187839
$ws.Cells.Item(17, 10).Value = 0.3290118330296043
$ws.Cells.Item(17, 11).Value = 0.2805100746741941
$ws.Cells.Item(17, 12).Value = 1.253652661228271
$ws.Cells.Item(17, 13).Value = 0.5290924020623156
$ws.Cells.Item(17, 14).Value = 1.017626570149982
$ws.Cells.Item(17, 15).Value = 0.5516169609903697
$ws.Cells.Item(17, 16).Value = 164.546368757167
$ws.Cells.Item(17, 17).Value = 263.2753105714912

$ws.Cells.Item(18, 1).Value = "model_21_9_8"
$ws.Cells.Item(18, 2).Value = 0.9579635622922461
$ws.Cells.Item(18, 3).Value = 0.8131640746735642
$ws.Cells.Item(18, 4).Value = 0.912253498774078
$ws.Cells.Item(18, 5).Value = 0.7381558855372499
$ws.Cells.Item(18, 6).Value = 0.8548567536465963
$ws.Cells.Item(18, 7).Value = 0.2810978525117478
$ws.Cells.Item(18, 8).Value = 1.249372692958207
$ws.Cells.Item(18, 9).Value = 0.2152090935514031
$ws.Cells.Item(18, 10).Value = 0.3289178396019615
$ws.Cells.Item(18, 11).Value = 0.2720634665766823
$ws.Cells.Item(18, 12).Value = 1.250039618575056
$ws.Cells.Item(18, 13).Value = 0.5301866204571253
$ws.Cells.Item(18, 14).Value = 1.017699552719054
$ws.Cells.Item(18, 15).Value = 0.5527577625275908
$ws.Cells.Item(18, 16).Value = 164.5381048813741
$ws.Cells.Item(18, 17).Value = 263.2670466956984

$ws.Cells.Item(19, 1).Value = "model_21_9_7"
$ws.Cells.Item(19, 2).Value = 0.9577572555628689
$ws.Cells.Item(19, 3).Value = 0.8128054601740308
$ws.Cells.Item(19, 4).Value = 0.9188225403618127
$ws.Cells.Item(19, 5).Value = 0.7380727209343219
$ws.Cells.Item(19, 6).Value = 0.8591265424330405
$ws.Cells.Item(19, 7).Value = 0.2824774265610484
$ws.Cells.Item(19, 8).Value = 1.251770749767855
$ws.Cells.Item(19, 9).Value = 0.1990977105806119
$ws.Cells.Item(19, 10).Value = 0.3290223075659733
$ws.Cells.Item(19, 11).Value = 0.2640599695626925
$ws.Cells.Item(19, 12).Value = 1.24647410754024
$ws.Cells.Item(19, 13).Value = 0.5314860549074155
$ws.Cells.Item(19, 14).Value = 1.017786418710371
$ws.Cells.Item(19, 15).Value = 0.5541125165926302
$ws.Cells.Item(19, 16).Value = 164.5283132751338
$ws.Cells.Item(19, 17).Value = 263.257255089458

$ws.Cells.Item(20, 1).Value = "model_21_9_6"
$ws.Cells.Item(20, 2).Value = 0.9575175072469765
$ws.Cells.Item(20, 3).Value = 0.8124178718058737
$ws.Cells.Item(20, 4).Value = 0.9251069869367774
$ws.Cells.Item(20, 5).Value = 0.737829166224698
$ws.Cells.Item(20, 6).Value = 0.8631564172502991
$ws.Cells.Item(20, 7).Value = 0.284080624653361
$ws.Cells.Item(20, 8).Value = 1.254362554970402
$ws.Cells.Item(20, 9).Value = 0.1836843319048266
$ws.Cells.Item(20, 10).Value = 0.3293282510051786
$ws.Cells.Item(20, 11).Value = 0.2565061788063254
$ws.Cells.Item(20, 12).Value = 1.242964975815851
$ws.Cells.Item(20, 13).Value = 0.532992143144119
$ws.Cells.Item(20, 14).Value = 1.017887365369694
$ws.Cells.Item(20, 15).Value = 0.5556827221236025
$ws.Cells.Item(20, 16).Value = 164.5169943829709
$ws.Cells.Item(20, 17).Value = 263.2459361972951

$ws.Cells.Item(21, 1).Value = "model_21_9_5"
$ws.Cells.Item(21, 2).Value = 0.9572437238055154
$ws.Cells.Item(21, 3).Value = 0.8120008469940204
$ws.Cells.Item(21, 4).Value = 0.9311030973630244
$ws.Cells.Item(21, 5).Value = 0.737424111146467
$ws.Cells.Item(21, 6).Value = 0.8669434130000799
$ws.Cells.Item(21, 7).Value = 0.2859114157870677
$ws.Cells.Item(21, 8).Value = 1.257151201807488
$ws.Cells.Item(21, 9).Value = 0.1689781331203956
$ws.Cells.Item(21, 10).Value = 0.3298370645850634
$ws.Cells.Item(21, 11).Value = 0.2494076522301187
$ws.Cells.Item(21, 12).Value = 1.239521044611037
$ws.Cells.Item(21, 13).Value = 0.5347068503274179
$ws.Cells.Item(21, 14).Value = 1.018002642608204
$ws.Cells.Item(21, 15).Value = 0.5574704279416276
$ws.Cells.Item(21, 16).Value = 164.5041465023007
$ws.Cells.Item(21, 17).Value = 263.2330883166249

$ws.Cells.Item(22, 1).Value = "model_21_9_4"
$ws.Cells.Item(22, 2).Value = 0.9569355508091286
$ws.Cells.Item(22, 3).Value = 0.8115540284364282
$ws.Cells.Item(22, 4).Value = 0.9368075540523553
$ws.Cells.Item(22, 5).Value = 0.7368559905152475
$ws.Cells.Item(22, 6).Value = 0.8704850256165187
$ws.Cells.Item(22, 7).Value = 0.2879721700329125
$ws.Cells.Item(22, 8).Value = 1.260139079559517
$ws.Cells.Item(22, 9).Value = 0.1549872510206833
$ws.Cells.Item(22, 10).Value = 0.3305507144260667
$ws.Cells.Item(22, 11).Value = 0.242769083575302
$ws.Cells.Item(22, 12).Value = 1.236137734211465
$ws.Cells.Item(22, 13).Value = 0.5366303849326018
$ws.Cells.Item(22, 14).Value = 1.018132399659314
$ws.Cells.Item(22, 15).Value = 0.5594758513972199
$ws.Cells.Item(22, 16).Value = 164.489782870691
$ws.Cells.Item(22, 17).Value = 263.2187246850153

$ws.Cells.Item(23, 1).Value = "model_21_9_3"
$ws.Cells.Item(23, 2).Value = 0.9565923058614174
$ws.Cells.Item(23, 3).Value = 0.8110769841618465
$ws.Cells.Item(23, 4).Value = 0.9422171220473886
$ws.Cells.Item(23, 5).Value = 0.7361209277364253
$ws.Cells.Item(23, 6).Value = 0.8737778978545616
$ws.Cells.Item(23, 7).Value = 0.2902674505787551
$ws.Cells.Item(23, 8).Value = 1.263329076820232
$ws.Cells.Item(23, 9).Value = 0.1417196197368061
$ws.Cells.Item(23, 10).Value = 0.3314740701473822
$ws.Cells.Item(23, 11).Value = 0.2365967658231227
$ws.Cells.Item(23, 12).Value = 1.232826036417028
$ws.Cells.Item(23, 13).Value = 0.5387647451149297
$ws.Cells.Item(23, 14).Value = 1.018276923847824
$ws.Cells.Item(23, 15).Value = 0.5617010757112439
$ws.Cells.Item(23, 16).Value = 164.4739050752635
$ws.Cells.Item(23, 17).Value = 263.2028468895877

$ws.Cells.Item(24, 1).Value = "model_21_9_2"
$ws.Cells.Item(24, 2).Value = 0.9562136826103204
$ws.Cells.Item(24, 3).Value = 0.8105693403078071
$ws.Cells.Item(24, 4).Value = 0.9473294027159709
$ws.Cells.Item(24, 5).Value = 0.7352200870038086
$ws.Cells.Item(24, 6).Value = 0.8768205909976691
$ws.Cells.Item(24, 7).Value = 0.2927993060022405
$ws.Cells.Item(24, 8).Value = 1.266723693609679
$ws.Cells.Item(24, 9).Value = 0.1291811222093295
$ws.Cells.Item(24, 10).Value = 0.3326056693364866
$ws.Cells.Item(24, 11).Value = 0.230893395772908
$ws.Cells.Item(24, 12).Value = 1.229590373951794
$ws.Cells.Item(24, 13).Value = 0.5411093290659851
$ws.Cells.Item(24, 14).Value = 1.018436344164076
$ws.Cells.Item(24, 15).Value = 0.5641454734550536
$ws.Cells.Item(24, 16).Value = 164.4565357342174
$ws.Cells.Item(24, 17).Value = 263.1854775485417

$ws.Cells.Item(25, 1).Value = "model_21_9_1"
$ws.Cells.Item(25, 2).Value = 0.9557992339479494
$ws.Cells.Item(25, 3).Value = 0.8100306605746124
$ws.Cells.Item(25, 4).Value = 0.9521422051873649
$ws.Cells.Item(25, 5).Value = 0.7341492727651543
$ws.Cells.Item(25, 6).Value = 0.8796104093816058
$ws.Cells.Item(25, 7).Value = 0.2955707261158759
$ws.Cells.Item(25, 8).Value = 1.270325847465944
$ws.Cells.Item(25, 9).Value = 0.117377131818375
$ws.Cells.Item(25, 10).Value = 0.3339507822740672
$ws.Cells.Item(25, 11).Value = 0.2256640263070694
$ws.Cells.Item(25, 12).Value = 1.226433542329786
$ws.Cells.Item(25, 13).Value = 0.5436641666653007
$ws.Cells.Item(25, 14).Value = 1.018610848864021
$ws.Cells.Item(25, 15).Value = 0.5668090757802151
$ws.Cells.Item(25, 16).Value = 164.4376942537928
$ws.Cells.Item(25, 17).Value = 263.166636068117

$ws.Cells.Item(26, 1).Value = "model_21_9_0"
$ws.Cells.Item(26, 2).Value = 0.9553484417697699
$ws.Cells.Item(26, 3).Value = 0.8094605760120754
$ws.Cells.Item(26, 4).Value = 0.9566533345946867
$ws.Cells.Item(26, 5).Value = 0.7329084114714396
$ws.Cells.Item(26, 6).Value = 0.8821459816569769
$ws.Cells.Item(26, 7).Value = 0.2985851754870696
$ws.Cells.Item(26, 8).Value = 1.274138005560626
$ws.Cells.Item(26, 9).Value = 0.1063130317451062
$ws.Cells.Item(26, 10).Value = 0.3355095013493908
$ws.Cells.Item(26, 11).Value = 0.2209112279487254
$ws.Cells.Item(26, 12).Value = 1.223685804479678
$ws.Cells.Item(26, 13).Value = 0.5464294789696742
$ws.Cells.Item(26, 14).Value = 1.018800656096939
$ws.Cells.Item(26, 15).Value = 0.5696921131543714
$ws.Cells.Item(26, 16).Value = 164.4174000836746
$ws.Cells.Item(26, 17).Value = 263.1463418979989
